# Update crypto price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $newValue
    $range.Style = "Normal"
}

Set-TextValue "D2" "26.316.61"
Set-TextValue "E2" "  +0.49%  "
Set-TextValue "D3" "1.596.48"
Set-TextValue "E3" "  +0.40%  "
Set-TextValue "E4" "  -0.03%  "
Set-TextValue "D5" "211.53"
Set-TextValue "E5" "  -0.12%  "
Set-TextValue "E6" "  -0.10%  "
Set-TextValue "E7" "  -0.03%  "
Set-TextValue "E8" "  +0.06%  "
Set-TextValue "E9" "  -0.09%  "
Set-TextValue "D10" "19.05"
Set-TextValue "E10" "  +0.16%  "
Set-TextValue "D11" "0.0854"
Set-TextValue "E11" "  +1.13%  "
Set-TextValue "D12" "1.821.04"
Set-TextValue "E12" "  +0.39%  "
Set-TextValue "D13" "1.598.99"
Set-TextValue "E13" "  +0.60%  "
Set-TextValue "D14" "3.99"
Set-TextValue "E14" "  -0.63%  "
Set-TextValue "D15" "0.504"
Set-TextValue "E15" "  -1.24%  "
Set-TextValue "D16" "63.47"
Set-TextValue "E16" "  -0.24%  "
Set-TextValue "D17" "26.305.24"
Set-TextValue "E17" "  +0.51%  "
Set-TextValue "D18" "230.36"
Set-TextValue "E18" "  +7.31%  "
Set-TextValue "E19" "  +3.86%  "
Set-TextValue "E20" "  -0.46%  "
Set-TextValue "E21" "  -0.01%  "
Set-TextValue "E22" "  -0.15%  "
Set-TextValue "E23" "  +2.46%  "
Set-TextValue "D24" "8.93"
Set-TextValue "E24" "  -0.93%  "
Set-TextValue "D25" "146.43"
Set-TextValue "E25" "  +1.10%  "
Set-TextValue "E26" "  -0.02%  "
Set-TextValue "D27" "6.96"
Set-TextValue "E27" "  +0.03%  "
Set-TextValue "E28" "  +0.04%  "
Set-TextValue "E29" "  +2.21%  "
Set-TextValue "D30" "0.0494"
Set-TextValue "E30" "  +0.10%  "
Set-TextValue "E31" "  +0.19%  "
Set-TextValue "D32" "1.505.24"
Set-TextValue "E32" "  +6.09%  "
Set-TextValue "E33" "  +0.96%  "
Set-TextValue "E34" "  -0.92%  "
Set-TextValue "E35" "  -0.36%  "
Set-TextValue "E36" "  +0.40%  "
Set-TextValue "E37" "  -3.33%  "
Set-TextValue "E38" "  -0.32%  "
Set-TextValue "E39" "  -1.03%  "
Set-TextValue "D40" "5.79"
Set-TextValue "E40" "  -1.52%  "
Set-TextValue "E41" "  +0.02%  "
Set-TextValue "B42" "WEMIXToken"
Set-TextValue "C42" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D42" "0.942"
Set-TextValue "E42" "  -2.18%  "
Set-TextValue "B43" "MXToken"
Set-TextValue "C43" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D43" "2.16"
Set-TextValue "E43" "  +1.29%  "
Set-TextValue "D44" "1.733.28"
Set-TextValue "E44" "  +0.45%  "
Set-TextValue "D45" "0.759"
Set-TextValue "E45" "  -0.52%  "
Set-TextValue "D46" "60.49"
Set-TextValue "E46" "  -0.90%  "
Set-TextValue "D47" "88.43"
Set-TextValue "E47" "  +1.54%  "
Set-TextValue "E48" "  -0.36%  "
Set-TextValue "E49" "  -0.18%  "
Set-TextValue "D50" "0.0956"
Set-TextValue "E50" "  -0.38%  "
Set-TextValue "E51" "  +0.06%  "
